# Updates crypto price (D) and volume-change (E) columns to match the
# latest scrape. D-column values that look numeric need NumberFormat
# forced to Text first so Excel doesn't renormalize the literal string
# (e.g. keep trailing zeros like '1.00' / '3.90' instead of '1' / '3.9').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.488.41'
$ws.Range('E2').Value = '  -0.63%  '

$ws.Range('D3').Value = '2.497.67'
$ws.Range('E3').Value = '  -0.75%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.29'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.14%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  -1.69%  '

$ws.Range('D9').Value = '2.495.70'
$ws.Range('E9').Value = '  -0.80%  '

$ws.Range('E10').Value = '  -0.55%  '

$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.356'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.85%  '

$ws.Range('E13').Value = '  +0.97%  '

$ws.Range('D14').Value = '2.955.82'
$ws.Range('E14').Value = '  -0.47%  '

$ws.Range('D15').Value = '69.444.20'
$ws.Range('E15').Value = '  -0.49%  '

$ws.Range('E16').Value = '  +0.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.60%  '

$ws.Range('D18').Value = '2.494.38'
$ws.Range('E18').Value = '  -0.85%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.38%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.06%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '347.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.61%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.61%  '

$ws.Range('E24').Value = '  -0.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.72%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.36%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.63%  '

$ws.Range('D28').Value = '2.628.69'
$ws.Range('E28').Value = '  -0.61%  '

$ws.Range('E29').Value = '  +0.31%  '

$ws.Range('D30').Value = '0.0₃0889'
$ws.Range('E30').Value = '  -2.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.02%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '457.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.86%  '

$ws.Range('E33').Value = '  -4.82%  '

$ws.Range('E34').Value = '  -1.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.115'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.55%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.14'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.03'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.37'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.38%  '

$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.316'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.96%  '

$ws.Range('E42').Value = '  -2.28%  '

$ws.Range('E43').Value = '  -0.99%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '

$ws.Range('E45').Value = '  -4.69%  '

$ws.Range('E46').Value = '  -6.71%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.92'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.75%  '

$ws.Range('E48').Value = '  -0.57%  '

$ws.Range('E49').Value = '  -2.50%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0731'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.575'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.56%  '
